# Fixed a bug in calcScatterProbability
# The fix reorders the data rows (A2:F25) on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(902,  1,  0,  0,  0,  0),
    @(201,  9, 30, 15, 45, 30),
    @(801,  3, 67, 65, 52, 45),
    @(101,  9, 30, 15, 60, 15),
    @(501,  9, 52, 30, 75, 45),
    @(1201, 2, 10, 10, 10, 10),
    @(1202, 2, 10, 10, 10, 10),
    @(901, 16, 15, 45, 60, 60),
    @(1001,18, 30, 75, 60, 72),
    @(601,  9, 60, 67, 60, 42),
    @(301,  6, 45, 30, 60, 45),
    @(401,  9, 48, 67, 75, 45),
    @(701,  3, 90, 45, 97, 15),
    @(1203, 3, 15, 15, 15, 15),
    @(502,  0,  4,  0,  0,  0),
    @(3,    0,  3,  3,  3,  3),
    @(1101, 0, 15, 30, 30,  0),
    @(802,  0,  4,  5,  4,  0),
    @(1,    0,  2,  2,  2,  2),
    @(2,    0,  2,  2,  2,  2),
    @(402,  0,  0,  4,  0,  0),
    @(602,  0,  0,  4,  0,  9),
    @(702,  0,  0,  0,  4,  0),
    @(1002, 0,  0,  0,  0,  9)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $rowValues.Length; $j++) {
        $ws.Cells.Item($rowNum, $j + 1).Value = $rowValues[$j]
    }
}
